# The sheet currently has two duplicated/derivative columns (F: "产销率",
# G: "销售量") that mirror column B / column E and are being dropped, and
# for every year the "B" (second) and "C" (third) quarterly rows had been
# recorded swapped relative to the source - this script fixes the row
# order and removes the now-redundant trailing columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (B-row, C-row) to swap for each year block of 4 rows
# (A,B,C,D quarters), starting at row 2.
$pairs = @(
    @(3,4), @(7,8), @(11,12), @(15,16), @(19,20), @(23,24), @(27,28),
    @(31,32), @(35,36), @(39,40), @(43,44), @(47,48), @(51,52), @(55,56),
    @(59,60), @(63,64), @(67,68)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("A$r1" + ":E$r1")
    $range2 = $ws.Range("A$r2" + ":E$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Drop the two trailing derived columns (F: product-sales ratio duplicate,
# G: sales-volume duplicate) now that the row order has been corrected.
$ws.Range("F:G").Delete()
